$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Balon" / NumProduct / Validation) is refreshed with a new scrape
# result: the NumProduct count becomes "3" and the Validation text is
# cleared out.
$ws.Range("B2").Value = "3"
$ws.Range("C2").Value = ""
